# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计" holding
#    the per-fund holdings detail for the new quarter.
# 2. Insert a new leading row into "总计" summarizing the new quarter and
#    bump the existing row index counters down by one.

$wb = $excel.ActiveWorkbook

# Writes $value into $range as literal TEXT, even when it looks like a
# number (e.g. "011815" or "4.64") - matching the source data, which
# stores these columns as strings, not numbers. Toggling to a text
# number-format before the write keeps Excel from auto-converting it to
# a numeric literal; ClearFormats() afterwards drops that temporary
# format again so no stray style sticks to the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# 1) Create the "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------
$q4_2021 = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4_2021)
$newSheet.Name = "2022-Q1"

Set-TextValue $newSheet.Range("B1") "基金代码"
Set-TextValue $newSheet.Range("C1") "基金名称"
Set-TextValue $newSheet.Range("D1") "基金规模"
Set-TextValue $newSheet.Range("E1") "股票总仓位"
Set-TextValue $newSheet.Range("F1") "仓位占比"
Set-TextValue $newSheet.Range("G1") "持有市值(亿元)"
Set-TextValue $newSheet.Range("H1") "仓位排名"

$rows = @(
    @(0, "011815", "恒越优势精选混合型发起式证券投资基金", "4.64", "92.44", "4.73", "0.2195", 3),
    @(1, "003318", "景顺长城中证500行业中性低波动指数", "13.99", "93.88", "1.22", "0.1707", 5),
    @(2, "013028", "恒越品质生活混合", "2.03", "92.89", "4.40", "0.0893", 4),
    @(3, "519677", "银河定投宝中证腾讯济安价值100A股指数", "2.74", "91.56", "1.28", "0.0351", 5),
    @(4, "512260", "华安中证500行业中性低波动ETF", "1.17", "96.94", "1.26", "0.0147", 5),
    @(5, "009658", "汇丰晋信中小盘低波动策略股票A", "0.98", "86.56", "1.30", "0.0127", 1),
    @(6, "007943", "富安达中证 500 指数增强", "0.21", "93.50", "1.40", "0.0029", 3),
    @(7, "009775", "汇丰晋信中小盘低波动策略股票C", "0.04", "86.56", "1.30", "0.0005", 1)
)

foreach ($r in $rows) {
    $rowIdx = 2 + $r[0]
    $newSheet.Range("A$rowIdx").Value = $r[0]
    Set-TextValue $newSheet.Range("B$rowIdx") $r[1]
    Set-TextValue $newSheet.Range("C$rowIdx") $r[2]
    Set-TextValue $newSheet.Range("D$rowIdx") $r[3]
    Set-TextValue $newSheet.Range("E$rowIdx") $r[4]
    Set-TextValue $newSheet.Range("F$rowIdx") $r[5]
    Set-TextValue $newSheet.Range("G$rowIdx") $r[6]
    $newSheet.Range("H$rowIdx").Value = $r[7]
}

# Copy header + index-column formatting (style "2") from the "2021-Q4"
# sheet so the new sheet matches the look of its siblings. Done AFTER
# the text writes above (which strip styles via ClearFormats) so the
# copied formatting is the one that sticks.
$q4_2021.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q4_2021.Range("A2").Copy()
$newSheet.Range("A2:A9").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Prepend the 2022-Q1 summary row into "总计"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows("2:2").Insert()

# Re-apply the "index column" / data styling that Insert() only partially
# carried over, so the new row matches the rest of the table exactly.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B3:D3").Copy()
$total.Range("B2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 8
$total.Range("D2").Value = 0.55

# Bump the index counters of the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
